# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 7
    $ws.Range("F4").Value = 947
    $ws.Range("F5").Value = 218
    $ws.Range("F6").Value = 431
}
